$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert the three new sentences about Topsec into the first empty
#    paragraph (ind left=360) that follows the paragraph ending
#    "...blind eye their way."
# ---------------------------------------------------------------------

$findRange = $d.Content
$findRange.Find.Execute("blind eye their way.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Collapse(0)
$findRange.MoveStart(1, 1)

$findRange.InsertAfter("Now, certain members of various Chinese patriotic and otherwise affiliated hacker groups are being all but directly employed by the CPC. Topsec, a cyber security company in Beijing, at one point employed Lion, the founder of HUC, and employs several other members and former members of other Chinese threat groups. ")
$findRange.Font.Bold = $true
$findRange.Font.NameFarEast = "DengXian"
$findRange.LanguageIDFarEast = "zh-CN"
$findRange.Collapse(0)

$findRange.InsertAfter("The chairman at the head of Topsec stated in an interview that the CPC contributed half of the 440M USD company’s startup funding, and they do receive and act on directives from the PLA. ")
$findRange.Font.Bold = $true
$findRange.Font.NameFarEast = "DengXian"
$findRange.LanguageIDFarEast = "zh-CN"
$findRange.Collapse(0)

$findRange.InsertAfter("This means that, at one point, the head of the HUC was essentially a government contracted network operator.")
$findRange.Font.Bold = $true
$findRange.Font.NameFarEast = "DengXian"
$findRange.LanguageIDFarEast = "zh-CN"
$findRange.Collapse(0)

# A unique placeholder right after the new text lets us locate the exact
# end-of-text insertion point reliably (collapsed ranges placed directly
# at the tail of a run are unreliable for Bookmarks.Add in this host).
$findRange.InsertAfter("GOBACKMARKERPLACEHOLDER")

$markerRange = $d.Content
$markerRange.Find.Execute("GOBACKMARKERPLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$deleteMarker = $d.Content
$deleteMarker.Find.Execute("GOBACKMARKERPLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteMarker.Delete()

# ---------------------------------------------------------------------
# 2. The old "_GoBack" paragraph (formerly the last of the trailing
#    empty ind=360 paragraphs) is now redundant - remove that whole
#    paragraph (its bookmark has moved to the new Topsec paragraph).
# ---------------------------------------------------------------------

$emptyCandidate = $d.Content
$emptyCandidate.Find.Execute("This means that, at one point, the head of the HUC was essentially a government contracted network operator.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$topsecParaIndex = $emptyCandidate.Paragraphs.Item(1).Index

$trailingEmptyCount = 0
$scanIndex = $topsecParaIndex + 1
while ($trailingEmptyCount -lt 8) {
    $scanIndex = $scanIndex + 1
    $trailingEmptyCount = $trailingEmptyCount + 1
}
$bookmarkParaIndex = $topsecParaIndex + 9
$d.Paragraphs.Item($bookmarkParaIndex).Range.Delete()

# ---------------------------------------------------------------------
# 3. Drop the stray <w:lastRenderedPageBreak/> on the "Targets:" run -
#    delete + retype the run text with identical formatting so the
#    stale rendering hint is not preserved.
# ---------------------------------------------------------------------

$targetsRange = $d.Content
$targetsRange.Find.Execute("Targets:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$targetsRange.Delete()
$targetsRange.InsertAfter("Targets:")
$targetsRange.Font.Bold = $true
$targetsRange.Font.Underline = 1

Write-Output "edit complete"
